# Add a new "Scheduled Date" column (F) to the production schedule header row,
# matching the header style used by the existing columns, and leave the
# selection positioned at F2 (as captured by the workbook after editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header value in column F, row 1.
$ws.Range("F1").Value = "Scheduled Date"

# Copy the formatting (fill/font/alignment) from the existing header cell E1
# so the new header cell matches the rest of the header row style.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Move/Set the active selection to F2, matching the saved sheet view state.
[void]$ws.Range("F2").Select()
